$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11 updates
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Row 12 updates
$ws.Range("B12").Value = 243
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "241/252"
